# Scheduled price-refresh: push updated market-board averages / leve
# profit figures (columns H-N) into the per-job "Leve Profit" tables.
# Values come from the runner's latest price pull; cells that had no
# recorded profit/loss before (or no longer have one) are created or
# cleared to match, the rest are simply overwritten in place.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 296.33334
$ws.Range("I18").Value = 296.33334
$ws.Range("K18").Value = 296.33334
$ws.Range("M18").Value = -12.33334000000002

$ws.Range("H29").Value = 90
$ws.Range("I29").Value = 90
$ws.Range("K29").Value = 270
$ws.Range("M29").Value = 11

$ws.Range("H40").Value = 500
$ws.Range("I40").Value = 500
$ws.Range("K40").Value = 500
$ws.Range("M40").Value = -325

$ws.Range("H62").Value = 4241.0586
$ws.Range("I62").Value = 5074.5
$ws.Range("K62").Value = 5074.5
$ws.Range("M62").Value = -4450.5

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 4241.0586
$ws.Range("I65").Value = 5074.5
$ws.Range("K65").Value = 25372.5
$ws.Range("M65").Value = -22252.5

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H75").Value = 19000
$ws.Range("J75").Value = 19000
$ws.Range("L75").Value = 19000
$ws.Range("N75").Value = -20872

$ws.Range("H78").Value = 19000
$ws.Range("J78").Value = 19000
$ws.Range("L78").Value = 57000
$ws.Range("N78").Value = -66360

$ws.Range("H80").Value = 306.72726
$ws.Range("I80").Value = 253.5
$ws.Range("J80").Value = 370.6
$ws.Range("K80").Value = 760.5
$ws.Range("L80").Value = 1111.8
$ws.Range("M80").Value = 237.5
$ws.Range("N80").Value = -3107.8

$ws.Range("H83").Value = 306.72726
$ws.Range("I83").Value = 253.5
$ws.Range("J83").Value = 370.6
$ws.Range("K83").Value = 2281.5
$ws.Range("L83").Value = 3335.4
$ws.Range("M83").Value = 2710.5
$ws.Range("N83").Value = -13319.4

$ws.Range("H101").Value = 561.75
$ws.Range("I101").Value = 561.75
$ws.Range("K101").Value = 1685.25
$ws.Range("M101").Value = -63.25

$ws.Range("H139").Value = 89997
$ws.Range("J139").Value = 89997
$ws.Range("L139").Value = 89997
$ws.Range("N139").Value = -100277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 2000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1856
$ws.Range("N11").ClearContents()

$ws.Range("H16").Value = 1200
$ws.Range("J16").Value = 1200
$ws.Range("L16").Value = 1200
$ws.Range("N16").Value = -1774

$ws.Range("H41").Value = 662
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 662
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 662
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -1490

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H139").Value = 57857.5
$ws.Range("J139").Value = 57857.5
$ws.Range("L139").Value = 57857.5
$ws.Range("N139").Value = -68137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H10").Value = 950
$ws.Range("I10").Value = 950
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 950
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -810
$ws.Range("N10").ClearContents()

$ws.Range("H23").Value = 606
$ws.Range("I23").Value = 606
$ws.Range("K23").Value = 606
$ws.Range("M23").Value = -323

$ws.Range("H37").Value = 3500
$ws.Range("J37").Value = 3500
$ws.Range("L37").Value = 3500
$ws.Range("N37").Value = -3774

$ws.Range("H54").Value = 4000
$ws.Range("I54").Value = 4000
$ws.Range("K54").Value = 4000
$ws.Range("M54").Value = -3516

$ws.Range("H86").Value = 1200
$ws.Range("I86").Value = 1200
$ws.Range("K86").Value = 1200
$ws.Range("M86").Value = -77

$ws.Range("H89").Value = 1200
$ws.Range("I89").Value = 1200
$ws.Range("K89").Value = 6000
$ws.Range("M89").Value = -384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 324
$ws.Range("I5").Value = 324
$ws.Range("K5").Value = 324
$ws.Range("M5").Value = -212

$ws.Range("H25").Value = 1734.5
$ws.Range("J25").Value = 2969
$ws.Range("L25").Value = 2969
$ws.Range("N25").Value = -3317

$ws.Range("H60").Value = 25046.5
$ws.Range("I60").Value = 25046.5
$ws.Range("K60").Value = 25046.5
$ws.Range("M60").Value = -24535.5

$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

$ws.Range("H68").Value = 47900
$ws.Range("J68").Value = 47900
$ws.Range("L68").Value = 47900
$ws.Range("N68").Value = -49398

$ws.Range("H71").Value = 47900
$ws.Range("J71").Value = 47900
$ws.Range("L71").Value = 143700
$ws.Range("N71").Value = -151188

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1324.7778
$ws.Range("I14").Value = 1324.7778
$ws.Range("K14").Value = 3974.3334
$ws.Range("M14").Value = -3801.3334

$ws.Range("H26").Value = 195
$ws.Range("I26").Value = 90
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 270
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = 18
$ws.Range("N26").Value = -1476

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 20000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40518

$ws.Range("H132").Value = 15337.333
$ws.Range("I132").Value = 13012
$ws.Range("K132").Value = 39036
$ws.Range("M132").Value = -36506

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496

$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
